# Reserva_salida1.xlsx - "Terminado hasta el primer informe"
#
# The "reserva_total.prn" sheet gets its "RESERVA PROGRAMADA EN EL PARQUE
# REGULANTE" block (rows 21-26) reshaped:
#   - rows 21-24 (HIDRO / TERMICA TG-CC / TERMICA TV / TOTAL) now reuse the
#     "... [MW]" labels already used earlier in the "POTENCIA OPERABLE"
#     block, instead of their own bare labels.
#   - rows 25-26 (RESERVA NUEVA / RESERVA TOTAL 2) become section-header
#     style rows: label text renamed with a " [MW]" suffix, merged A:E like
#     the other section headers (row 20, row 8, row 3, ...), bordered +
#     centered, and their value relocated from column D into column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reserva_total.prn")

# --- rows 21-24: relabel to reuse the "... [MW]" shared strings ---
$ws.Range("A21").Value = "HIDRO [MW]"
$ws.Range("A22").Value = "TÉRMICA TG-CC [MW]"
$ws.Range("A23").Value = "TÉRMICA TV [MW]"
$ws.Range("A24").Value = "TOTAL [MW]"

# --- row 25: RESERVA NUEVA -> RESERVA NUEVA [MW], header styling ---
$val25 = $ws.Range("D25").Value2
$ws.Range("D25").ClearContents()
$ws.Range("A25").Value = "RESERVA NUEVA [MW]"
$ws.Range("F25").Value = $val25
$ws.Range("A25:E25").Merge()
$ws.Range("A25:E25").HorizontalAlignment = -4108
$ws.Range("A25:E25").Borders.LineStyle = 1
$ws.Range("A25:E25").Borders.Weight = 2
$ws.Range("F25").HorizontalAlignment = -4108
$ws.Range("F25").Borders.LineStyle = 1
$ws.Range("F25").Borders.Weight = 2

# --- row 26: RESERVA TOTAL 2 -> RESERVA TOTAL 2 [MW], header styling ---
$val26 = $ws.Range("D26").Value2
$ws.Range("D26").ClearContents()
$ws.Range("A26").Value = "RESERVA TOTAL 2 [MW]"
$ws.Range("F26").Value = $val26
$ws.Range("A26:E26").Merge()
$ws.Range("A26:E26").HorizontalAlignment = -4108
$ws.Range("A26:E26").Borders.LineStyle = 1
$ws.Range("A26:E26").Borders.Weight = 2
$ws.Range("F26").HorizontalAlignment = -4108
$ws.Range("F26").Borders.LineStyle = 1
$ws.Range("F26").Borders.Weight = 2
